$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correlation-matrix values (rows 2-13, columns B-G) updated per the
# revised "macro_corr_educ_gr" abstract (learning-model correlations).
$values = @{
    'B2' = '-0.01'
    'C2' = '-0.05'
    'D2' = '-0.01'
    'E2' = '-0.04'
    'F2' = '-0.05'
    'G2' = '0.07'
    'B3' = '0.06'
    'C3' = '-0.03'
    'D3' = '0.03'
    'E3' = '0.12'
    'F3' = '-0.04'
    'G3' = '-0.16'
    'B4' = '-0.01'
    'C4' = '0.02'
    'D4' = '0.13'
    'E4' = '0.14'
    'F4' = '-0.04'
    'G4' = '-0.11'
    'B5' = '-0.0'
    'C5' = '-0.14'
    'D5' = '0.04'
    'E5' = '0.08'
    'F5' = '-0.12'
    'G5' = '-0.16'
    'B6' = '-0.11'
    'C6' = '-0.26**'
    'D6' = '-0.15'
    'E6' = '-0.11'
    'F6' = '0.03'
    'G6' = '0.04'
    'B7' = '-0.05'
    'C7' = '-0.27**'
    'D7' = '-0.12'
    'E7' = '-0.27**'
    'F7' = '-0.03'
    'G7' = '0.06'
    'B8' = '-0.08'
    'C8' = '-0.21*'
    'D8' = '-0.17'
    'E8' = '-0.11'
    'F8' = '-0.16'
    'G8' = '0.02'
    'B9' = '-0.07'
    'C9' = '-0.08'
    'D9' = '-0.08'
    'E9' = '-0.08'
    'F9' = '-0.21*'
    'G9' = '-0.07'
    'B10' = '-0.23*'
    'C10' = '0.07'
    'D10' = '-0.18'
    'E10' = '0.14'
    'F10' = '-0.14'
    'G10' = '-0.06'
    'B11' = '-0.28**'
    'C11' = '-0.01'
    'D11' = '-0.15'
    'E11' = '-0.03'
    'F11' = '0.01'
    'G11' = '0.02'
    'B12' = '-0.26**'
    'C12' = '-0.11'
    'D12' = '-0.08'
    'E12' = '-0.03'
    'F12' = '-0.08'
    'G12' = '-0.01'
    'B13' = '0.04'
    'C13' = '-0.16'
    'D13' = '0.07'
    'E13' = '-0.1'
    'F13' = '-0.19'
    'G13' = '0.07'
}

foreach ($addr in $values.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so values like "-0.0" / "-0.23*" keep their
    # literal representation instead of being coerced to numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $values[$addr]
    $cell.Style = "Normal"
}
